$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.524.40'
$ws.Range("E2").Value = '  +1.89%  '
$ws.Range("D3").Value = '2.966.50'
$ws.Range("E3").Value = '  +2.17%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.36'
$ws.Range("E5").Value = '  +1.63%  '
$ws.Range("E6").Value = '  +10.59%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.12%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.516'
$ws.Range("E8").Value = '  +4.02%  '
$ws.Range("D9").Value = '2.958.01'
$ws.Range("E9").Value = '  +2.21%  '
$ws.Range("E10").Value = '  +5.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4.82'
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.454'
$ws.Range("E12").Value = '  +3.82%  '
$ws.Range("E13").Value = '  +6.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.29'
$ws.Range("E14").Value = '  +2.87%  '
$ws.Range("E15").Value = '  +2.83%  '
$ws.Range("D16").Value = '3.454.76'
$ws.Range("E16").Value = '  +2.22%  '
$ws.Range("E17").Value = '  +5.11%  '
$ws.Range("D18").Value = '2.961.12'
$ws.Range("E18").Value = '  +2.39%  '
$ws.Range("D19").Value = '58.520.76'
$ws.Range("E19").Value = '  +2.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '419.35'
$ws.Range("E20").Value = '  +3.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.42'
$ws.Range("E21").Value = '  +4.01%  '
$ws.Range("E22").Value = '  +5.70%  '
$ws.Range("E23").Value = '  +3.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.26'
$ws.Range("E24").Value = '  +3.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '79.87'
$ws.Range("E25").Value = '  +3.71%  '
$ws.Range("E26").Value = '  -0.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.37%  '
$ws.Range("E28").Value = '  +8.75%  '
$ws.Range("E29").Value = '  +2.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.67'
$ws.Range("E30").Value = '  +6.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '25.47'
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.01'
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.996'
$ws.Range("E34").Value = '  +9.95%  '
$ws.Range("D35").Value = '0.0₃0745'
$ws.Range("E35").Value = '  +20.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.70'
$ws.Range("E36").Value = '  +5.37%  '
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("E38").Value = '  +1.48%  '
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.72'
$ws.Range("E40").Value = '  +12.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '395.27'
$ws.Range("E41").Value = '  +9.77%  '
$ws.Range("D42").Value = '2.736.16'
$ws.Range("E42").Value = '  +4.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0345'
$ws.Range("E43").Value = '  +0.89%  '
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("E45").Value = '  +0.05%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '124.68'
$ws.Range("E46").Value = '  +4.31%  '
$ws.Range("E47").Value = '  +5.39%  '
$ws.Range("E48").Value = '  +2.60%  '
$ws.Range("E49").Value = '  +1.88%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.74'
$ws.Range("E50").Value = '  +18.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.07'
$ws.Range("E51").Value = '  +0.87%  '
